$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.128.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.532.57"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.46%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.01"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.25"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.530.81"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.43%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.18%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.52%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.139.37"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000207"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.43"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.537.08"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.50%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.252.30"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.09%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.01"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.45%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.22"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.40%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.05"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "425.75"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.51%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.88"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.679.49"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.49%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.45%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.06"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.07%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.01%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.51"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.82%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.34"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.526.76"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.37%  "

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.76"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.87"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.62"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "171.30"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.33%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.12%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.892"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.31%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -9.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.31"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.23%  "

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.22"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -8.71%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.00"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -7.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.42"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.79%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.952"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.64%  "
